$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.205.11"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "1.603.28"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.40"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.31"
$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").Value = "1.826.00"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "1.603.97"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.512"
$ws.Range("E15").Value = "  -0.01%  "

$ws.Range("D16").Value = "26.179.74"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.83"
$ws.Range("E17").Value = "  +2.12%  "

$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "200.52"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("E22").Value = "  -0.65%  "

$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.24"
$ws.Range("E25").Value = "  +1.79%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("E27").Value = "  -2.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.18"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("E29").Value = "  +1.64%  "

$ws.Range("E30").Value = "  +3.42%  "

$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("E32").Value = "  +1.95%  "

$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("D36").Value = "1.161.71"
$ws.Range("E36").Value = "  +4.43%  "

$ws.Range("E37").Value = "  +3.47%  "

$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("E39").Value = "  -0.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.786"
$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("E43").Value = "  +3.74%  "

$ws.Range("D44").Value = "1.738.37"
$ws.Range("E44").Value = "  +0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.38"
$ws.Range("E45").Value = "  -1.70%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0106"
$ws.Range("E46").Value = "  +14.35%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.52"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("E51").Value = "  +0.10%  "
